# Updated symbol list on Sat Dec 31 06:57:19 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Price-only updates (column D) ----
# A leading apostrophe keeps these numeric-looking values stored as literal
# text (matching the source workbook's inlineStr cells) instead of being
# coerced into numbers / scientific notation.
$priceUpdates = [ordered]@{
    "D2"  = "245.06"
    "D3"  = "25.54"
    "D4"  = "5.113"
    "D5"  = "0.05584"
    "D6"  = "6.473"
    "D7"  = "3.015"
    "D9"  = "0.8404"
    "D22" = "3.742"
    "D23" = "0.04690"
    "D25" = "0.001246"
    "D26" = "0.004620"
    "D27" = "0.00009702"
    "D45" = "0.00005322"
    "D48" = "0.002123"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $priceUpdates[$addr]
}

# ---- Label-only updates (column E) ----
$ws.Range("E27").Value = "26NitroExNTXBestin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# ---- Full row re-rank for rows 11-20 (Coin / Link / Price / Volume) ----
$rowData = [ordered]@{
    11 = @{ B = "MandalaExchangeToken";               C = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx";               D = "0.06956";   E = "10MandalaExchangeTokenMDX" }
    12 = @{ B = "LiechtensteinCryptoassetsExchange";  C = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx";        D = "0.03229";   E = "11LiechtensteinCryptoassetsExchangeLCX" }
    13 = @{ B = "BitrueCoin";                          C = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr";                              D = "0.02878";   E = "12BitrueCoinBTR" }
    14 = @{ B = "BitMartToken";                        C = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx";                        D = "0.09390";   E = "13BitMartTokenBMX" }
    15 = @{ B = "BitForexToken";                       C = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf";                        D = "0.001514";  E = "14BitForexTokenBF" }
    16 = @{ B = "One";                                 C = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one";                                  D = "0.0005961"; E = "15OneONE" }
    17 = @{ B = "TigerCash";                           C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";                                D = "0.006183";  E = "16TigerCashTCH" }
    18 = @{ B = "LEO";                                 C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";                                  D = "3.524";     E = "17LEOLEO" }
    19 = @{ B = "BTSEToken";                           C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";                           D = "2.022";     E = "18BTSETokenBTSE" }
    20 = @{ B = "BitpandaEcosystemToken";              C = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best";                  D = "0.3179";    E = "19BitpandaEcosystemTokenBEST" }
}

foreach ($r in $rowData.Keys) {
    $row = $rowData[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = "'" + $row.D
    $ws.Range("E$r").Value = $row.E
}

# ---- Full row swap for rows 41-42 (Coin / Link / Price / Volume) ----
$swapRows = [ordered]@{
    41 = @{ B = "BKEXToken"; C = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk";         D = "0.1354";   E = "40BKEXTokenBKK" }
    42 = @{ B = "KickToken"; C = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick";    D = "0.006161"; E = "41KickTokenKICK" }
}

foreach ($r in $swapRows.Keys) {
    $row = $swapRows[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = "'" + $row.D
    $ws.Range("E$r").Value = $row.E
}
